$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 296, shifting existing rows 296-370 down to 297-371
$ws.Rows.Item(296).Insert()

# Populate the new row 296 with the new data record
$ws.Cells.Item(296, 1).Value = 11
$ws.Cells.Item(296, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(296, 3).Value = "Bíobío"
$ws.Cells.Item(296, 4).Value = 45275
$ws.Cells.Item(296, 5).Value = 8
$ws.Cells.Item(296, 6).Value = 100112040
$ws.Cells.Item(296, 7).Value = "Cilantro"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Primera"
$ws.Cells.Item(296, 10).Value = 80
$ws.Cells.Item(296, 11).Value = 6000
$ws.Cells.Item(296, 12).Value = 6000
$ws.Cells.Item(296, 13).Value = 6000
$ws.Cells.Item(296, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(296, 15).Value = "Región Metropolitana"
$ws.Cells.Item(296, 16).Value = 167
$ws.Cells.Item(296, 17).Value = 36
$ws.Cells.Item(296, 18).Value = "Hortaliza"
